$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell D6 ("预制茶销售量（单位）" column, row for 2024-01-15) mistakenly holds
# the text "4:36" (entered/interpreted as a time). Correct it to the plain
# numeric value 436 — this also drops the now-unused "4:36" shared string
# and Excel naturally renumbers the remaining shared-string table.
$ws.Range("D6").Value = 436
